# Update "想去人数" (F column) counts, cancel notice for 新余·原神&崩铁&崩坏only (C5),
# and mark its ticket price (G5) as "不可售" (not for sale) on both the
# "展览" and "全部类型" worksheets, which carry duplicate data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new F-column (想去人数) value
$fUpdates = @{
    3  = 27
    4  = 124
    5  = 115
    6  = 478
    9  = 605
    13 = 391
    15 = 107
    16 = 20
    18 = 55
    19 = 106
    20 = 1030
    22 = 312
    23 = 347
    30 = 269
    31 = 288
    32 = 1648
    36 = 595
    38 = 3820
    40 = 445
    42 = 956
    46 = 85
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Range("F$row").Value = $fUpdates[$row]
    }

    # Row 5 special-case updates: event cancelled
    $ws.Range("C5").Value = "新余·原神&崩铁&崩坏only（取消）"
    $ws.Range("G5").Value = "不可售"
}
